$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

$ws.Range("T1").Value = "Employee Type(Daily/Monthly/Fixed)"

$ws.Range("G1:G1048576").Validation.Delete()
$ws.Range("G2:G1048576").Validation.Add(3, 1, 1, """Male,Female""")

$ws.Range("I1:I1048576").Validation.Delete()
$ws.Range("I2:I1048576").Validation.Add(3, 1, 1, """Single,Married""")

$ws.Range("S1:S1048576").Validation.Delete()
$ws.Range("S2:S1048576").Validation.Add(3, 1, 1, """WEEKLY,SEMI-MONTHLY""")

$ws.Range("T1:T1048576").Validation.Delete()
$ws.Range("T2:T1048576").Validation.Add(3, 1, 1, """Daily,Monthly,Fixed""")

# Widen column T to fit new, longer header text (stored width 36)
$ws.Columns("T").ColumnWidth = 35.166666666666664

# Leave a trace of a (later cleared) filter on the Employee Type column,
# matching the hidden _FilterDatabase defined name Excel leaves behind
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=" + $ws.Name + "!T`$1:T`$1")
$fd.RefersTo = "=" + $ws.Name + "!`$T`$1:`$T`$1"
$fd.Visible = $false

# Move the active selection to the Employee Type column
$ws.Range("T2").Select()

# Re-protect the worksheet (content was protected before the edit)
$ws.Protect("CD03")

